$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "₹ 19,978"
$ws.Range("D2").Value = "₹ 127,861"
$ws.Range("C3").Value = "₹ 203,323"
$ws.Range("D3").Value = "₹ 1,283,575"
$ws.Range("C4").Value = "₹ 19,625"
$ws.Range("D4").Value = "₹ 96,160"
